# Add a new "UK" worksheet (FC market test case) as the first sheet,
# cloned from the existing "Italy" sheet template, then update its
# market name / ticket reference cells.

$wb = $excel.ActiveWorkbook

$italy = $wb.Worksheets.Item("Italy")

# Copy "Italy" and place the copy immediately before it -> becomes sheet #1.
$italy.Copy($italy)

$uk = $wb.Worksheets.Item(1)
$uk.Name = "UK"

# Update the market name and ticket reference for the new UK sheet.
$uk.Range("B2").Value = "UK Market"
$uk.Range("B4").Value = "NGC-3003/T1251/T1260"

# The previously active sheet (Hungary) is no longer the selected tab;
# its selection collapses back to the whole sheet.
$hungary = $wb.Worksheets.Item("Hungary")
$hungary.Activate()
$hungary.Cells.Select()

# UK becomes the new active tab, with A8 selected (mirrors the other
# sheets' "first data row under the notes" selection pattern).
$uk.Activate()
$uk.Range("A8").Select()
